$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet used to list "already registered" duplicate-email messages (A2:A5).
# It now reports a validation error ("name" is required) repeated for every
# unregistered student row, and the list has grown from 4 rows to 13 rows
# (A2:A14), which pushes the sheet's used range/dimension to A1:A14.
$message = '"name" is required for the student undefined'

for ($row = 2; $row -le 14; $row++) {
    $ws.Range("A$row").Value = $message
}
